# edit.ps1 - apply the graph_llm.pptx diff via PowerPoint COM-interop
$p = $ppt.ActivePresentation
$nl = [char]13

# ---------------------------------------------------------------
# 1) Slide 12 ("Moro Takeouts" / RAG slide) becomes the new
#    "Applications" slide with the taxonomy bullet list.
# ---------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = 'Applications'

$body12 = $s12.Shapes.Item(2).TextFrame.TextRange
$body12.Text = 'GraphText: Graph reasoning in text space - NIPS 2025:' + $nl + 'GraphText translates graphs to natural language. GRAPHTEXT derives a graph-syntax tree for each graph that encapsulates both the node attributes and inter-node relationships. Traversal of the tree yields a graph text sequence, which is then processed by an LLM to treat graph tasks as text generation tasks. Notably, GRAPHTEXT offers multiple advantages. It introduces training-free graph reasoning: even without training on graph data, GRAPHTEXT with ChatGPT can achieve on par with, or even surpassing, the performance of supervised-trained graph neural networks through in-context learning (ICL). Furthermore, GRAPHTEXT paves the way for interactive graph reasoning, allowing both humans and LLMs to communicate with the model seamlessly using natural language.' + $nl + 'Temporal Knowledge Graph Forecasting: predicting future facts occur at time tn based on given historical facts occur at time t with t < tn' + $nl + 'Knowledge base question answering:' + $nl + 'FlexKQBA: given a KG, create new data to train QA algorithms: given a query, parse them in natural language' + $nl + 'Knowledge Graph Question Generation (KGQG): generate questions based on graph.' + $nl + 'Knowledge graph completion:' + $nl + 'Class-Imbalanced Graph Learning: imbalance not only in labels (e.g. more fraud users than legit in a social network graph), but also in graph connectivity (more :marriedTo than :friendOf).' + $nl + 'LLM used to generate synthetic data to represent the unbalanced class;'

# second-level bullets (lvl="1" in OOXML == IndentLevel 2 in the OM)
$body12.Paragraphs(2,1).IndentLevel = 2
$body12.Paragraphs(5,1).IndentLevel = 2
$body12.Paragraphs(6,1).IndentLevel = 2
$body12.Paragraphs(9,1).IndentLevel = 2

# ---------------------------------------------------------------
# 2) Insert a brand-new slide at position 13 that carries the
#    original "Moro Takeouts" / RAG content (Title+Content layout).
# ---------------------------------------------------------------
$sNew = $p.Slides.Add(13, 2)
$sNew.Shapes.Item(1).TextFrame.TextRange.Text = 'Moro Takeouts'
$sNew.Shapes.Item(2).TextFrame.TextRange.Text = 'Use RAG only if the question lies beyond the typical training data, if the knowledge we’re asking is not embedded in parameters, otherwise performances might get worse'

# ---------------------------------------------------------------
# 3) The old slide 13 ("LLMs on Graphs" taxonomy slide) is now
#    slide 14; give it a second content placeholder with the
#    two new Italian notes.
# ---------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$srcPlaceholder = $s14.Shapes.Item(2)
$dup = $srcPlaceholder.Duplicate()
$newShape = $dup.Item(1)
$newShape.Name = "Content Placeholder 2"
$newShape.TextFrame.TextRange.Text = 'Una volta che embedd oe faccio rag e trovo qualcosa per similarità, lo passo comunque al LLM sotto forma testuale (?)' + $nl + 'se uso node2vec e embeddo un nodo, come glielo passo (e.g. text? json?)'

Write-Host "Slides total: " $p.Slides.Count
